# Applies the scheduled-runner price/profit updates to the Excalibur Profits sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
# Columns: H=currentAveragePrice, I=currentAveragePriceNQ, J=currentAveragePriceHQ,
#          K=LevePriceNQ, L=LevePriceHQ, M=LeveProfitNQ, N=LeveProfitHQ
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 9
$ws.Cells.Item(9, 8).Value = 442.3
$ws.Cells.Item(9, 9).Value = 419.22223
$ws.Cells.Item(9, 10).Value = 650
$ws.Cells.Item(9, 11).Value = 419.22223
$ws.Cells.Item(9, 12).Value = 650
$ws.Cells.Item(9, 13).Value = -250.22223
$ws.Cells.Item(9, 14).Value = -988

# Row 53
$ws.Cells.Item(53, 8).Value = 259.3158
$ws.Cells.Item(53, 9).Value = 250.07692
$ws.Cells.Item(53, 11).Value = 250.07692
$ws.Cells.Item(53, 13).Value = 386.92308

# Row 64
$ws.Cells.Item(64, 8).Value = 8539.583000000001
$ws.Cells.Item(64, 9).Value = 4831.6665
$ws.Cells.Item(64, 10).Value = 9775.556
$ws.Cells.Item(64, 11).Value = 4831.6665
$ws.Cells.Item(64, 12).Value = 9775.556
$ws.Cells.Item(64, 13).Value = -4583.6665
$ws.Cells.Item(64, 14).Value = -10271.556

# Row 67
$ws.Cells.Item(67, 8).Value = 8539.583000000001
$ws.Cells.Item(67, 9).Value = 4831.6665
$ws.Cells.Item(67, 10).Value = 9775.556
$ws.Cells.Item(67, 11).Value = 4831.6665
$ws.Cells.Item(67, 12).Value = 9775.556
$ws.Cells.Item(67, 13).Value = -3973.6665
$ws.Cells.Item(67, 14).Value = -11491.556

# Row 69
$ws.Cells.Item(69, 8).Value = 9085.556
$ws.Cells.Item(69, 10).Value = 9427.5
$ws.Cells.Item(69, 12).Value = 28282.5
$ws.Cells.Item(69, 14).Value = -30030.5

# Row 72
$ws.Cells.Item(72, 8).Value = 9085.556
$ws.Cells.Item(72, 10).Value = 9427.5
$ws.Cells.Item(72, 12).Value = 84847.5
$ws.Cells.Item(72, 14).Value = -93583.5

# Row 87
$ws.Cells.Item(87, 8).Value = 97492.25
$ws.Cells.Item(87, 10).Value = 99989.664
$ws.Cells.Item(87, 12).Value = 99989.664
$ws.Cells.Item(87, 14).Value = -102485.664

# Row 90
$ws.Cells.Item(90, 8).Value = 97492.25
$ws.Cells.Item(90, 10).Value = 99989.664
$ws.Cells.Item(90, 12).Value = 299968.992
$ws.Cells.Item(90, 14).Value = -312448.992

# Row 106
$ws.Cells.Item(106, 8).Value = 2944.9412
$ws.Cells.Item(106, 9).Value = 2826.6428
$ws.Cells.Item(106, 10).Value = 3497
$ws.Cells.Item(106, 11).Value = 2826.6428
$ws.Cells.Item(106, 12).Value = 3497
$ws.Cells.Item(106, 13).Value = -2195.6428
$ws.Cells.Item(106, 14).Value = -4759

# Row 112
$ws.Cells.Item(112, 8).Value = 9373.906000000001
$ws.Cells.Item(112, 10).Value = 10654.214
$ws.Cells.Item(112, 12).Value = 31962.642
$ws.Cells.Item(112, 14).Value = -34178.642

# Row 138
$ws.Cells.Item(138, 8).Value = 2191.8386
$ws.Cells.Item(138, 10).Value = 2918.7334
$ws.Cells.Item(138, 12).Value = 8756.200199999999
$ws.Cells.Item(138, 14).Value = -19036.2002

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Cells.Item(32, 8).Value = 5467132
$ws.Cells.Item(32, 9).Value = 5749894.5
$ws.Cells.Item(32, 11).Value = 5749894.5
$ws.Cells.Item(32, 13).Value = -5749607.5

# Row 61
$ws.Cells.Item(61, 8).Value = 11146767
$ws.Cells.Item(61, 10).Value = 6969
$ws.Cells.Item(61, 12).Value = 6969
$ws.Cells.Item(61, 14).Value = -7393

# Row 74
$ws.Cells.Item(74, 8).Value = 11823.7
$ws.Cells.Item(74, 9).Value = 12770.818
$ws.Cells.Item(74, 10).Value = 10666.111
$ws.Cells.Item(74, 11).Value = 12770.818
$ws.Cells.Item(74, 12).Value = 10666.111
$ws.Cells.Item(74, 13).Value = -11896.818
$ws.Cells.Item(74, 14).Value = -12414.111

# Row 77
$ws.Cells.Item(77, 8).Value = 11823.7
$ws.Cells.Item(77, 9).Value = 12770.818
$ws.Cells.Item(77, 10).Value = 10666.111
$ws.Cells.Item(77, 11).Value = 63854.09
$ws.Cells.Item(77, 12).Value = 53330.55500000001
$ws.Cells.Item(77, 13).Value = -59486.09
$ws.Cells.Item(77, 14).Value = -62066.55500000001

# Row 136
$ws.Cells.Item(136, 8).Value = 11146767
$ws.Cells.Item(136, 10).Value = 6969
$ws.Cells.Item(136, 12).Value = 20907
$ws.Cells.Item(136, 14).Value = -26007

$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Cells.Item(20, 8).Value = 1270.5454
$ws.Cells.Item(20, 9).Value = 924.53845
$ws.Cells.Item(20, 10).Value = 1770.3334
$ws.Cells.Item(20, 11).Value = 924.53845
$ws.Cells.Item(20, 12).Value = 1770.3334
$ws.Cells.Item(20, 13).Value = -677.53845
$ws.Cells.Item(20, 14).Value = -2264.3334

$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Cells.Item(16, 8).Value = 400
$ws.Cells.Item(16, 9).Value = 400
$ws.Cells.Item(16, 11).Value = 400
$ws.Cells.Item(16, 13).Value = -113

# Row 31
$ws.Cells.Item(31, 8).Value = 119112.2
$ws.Cells.Item(31, 9).Value = 522085.66
$ws.Cells.Item(31, 11).Value = 522085.66
$ws.Cells.Item(31, 13).Value = -521790.66

# Row 34
$ws.Cells.Item(34, 8).Value = 119112.2
$ws.Cells.Item(34, 9).Value = 522085.66
$ws.Cells.Item(34, 11).Value = 522085.66
$ws.Cells.Item(34, 13).Value = -521883.66

# Row 58
$ws.Cells.Item(58, 8).Value = 652544.3
$ws.Cells.Item(58, 9).Value = 1123657.9
$ws.Cells.Item(58, 10).Value = 4763.125
$ws.Cells.Item(58, 11).Value = 1123657.9
$ws.Cells.Item(58, 12).Value = 4763.125
$ws.Cells.Item(58, 13).Value = -1123454.9
$ws.Cells.Item(58, 14).Value = -5169.125

# Row 62
$ws.Cells.Item(62, 8).Value = 4921
$ws.Cells.Item(62, 9).Value = 2401.25
$ws.Cells.Item(62, 10).Value = 15000
$ws.Cells.Item(62, 11).Value = 2401.25
$ws.Cells.Item(62, 12).Value = 15000
$ws.Cells.Item(62, 13).Value = -1777.25
$ws.Cells.Item(62, 14).Value = -16248

# Row 65
$ws.Cells.Item(65, 8).Value = 4921
$ws.Cells.Item(65, 9).Value = 2401.25
$ws.Cells.Item(65, 10).Value = 15000
$ws.Cells.Item(65, 11).Value = 12006.25
$ws.Cells.Item(65, 12).Value = 75000
$ws.Cells.Item(65, 13).Value = -8886.25
$ws.Cells.Item(65, 14).Value = -81240

# Row 113
$ws.Cells.Item(113, 8).Value = 400
$ws.Cells.Item(113, 9).Value = 400
$ws.Cells.Item(113, 11).Value = 400
$ws.Cells.Item(113, 13).Value = 1770

# Row 134
$ws.Cells.Item(134, 8).Value = 11712.23
$ws.Cells.Item(134, 9).Value = 11712.23
$ws.Cells.Item(134, 11).Value = 35136.69
$ws.Cells.Item(134, 13).Value = -32601.69

# Row 136
$ws.Cells.Item(136, 8).Value = 652544.3
$ws.Cells.Item(136, 9).Value = 1123657.9
$ws.Cells.Item(136, 10).Value = 4763.125
$ws.Cells.Item(136, 11).Value = 3370973.7
$ws.Cells.Item(136, 12).Value = 14289.375
$ws.Cells.Item(136, 13).Value = -3368423.7
$ws.Cells.Item(136, 14).Value = -19389.375

$ws = $wb.Worksheets.Item("CUL")
# Row 58
$ws.Cells.Item(58, 8).Value = 740
$ws.Cells.Item(58, 10).Value = 790
$ws.Cells.Item(58, 12).Value = 2370
$ws.Cells.Item(58, 14).Value = -2626

$ws = $wb.Worksheets.Item("GSM")
# Row 26
$ws.Cells.Item(26, 8).Value = 40999
$ws.Cells.Item(26, 10).Value = 40999
$ws.Cells.Item(26, 12).Value = 40999
$ws.Cells.Item(26, 14).Value = -41559

# Row 50
$ws.Cells.Item(50, 8).Value = 40999
$ws.Cells.Item(50, 10).Value = 40999
$ws.Cells.Item(50, 12).Value = 40999
$ws.Cells.Item(50, 14).Value = -41995

# Row 52
$ws.Cells.Item(52, 8).Value = 0
$ws.Cells.Item(52, 10).Value = 0
$ws.Cells.Item(52, 12).Value = 0
$ws.Cells.Item(52, 14).ClearContents()

# Row 113
$ws.Cells.Item(113, 8).Value = 2171.6316
$ws.Cells.Item(113, 10).Value = 2826.2222
$ws.Cells.Item(113, 12).Value = 2826.2222
$ws.Cells.Item(113, 14).Value = -7166.2222

$ws = $wb.Worksheets.Item("LTW")
# Row 61
$ws.Cells.Item(61, 8).Value = 1668.909
$ws.Cells.Item(61, 9).Value = 1401.125
$ws.Cells.Item(61, 11).Value = 1401.125
$ws.Cells.Item(61, 13).Value = -1199.125

# Row 74
$ws.Cells.Item(74, 8).Value = 54409.285
$ws.Cells.Item(74, 9).Value = 19999
$ws.Cells.Item(74, 10).Value = 80217
$ws.Cells.Item(74, 11).Value = 19999
$ws.Cells.Item(74, 12).Value = 80217
$ws.Cells.Item(74, 13).Value = -19001
$ws.Cells.Item(74, 14).Value = -82213

# Row 77
$ws.Cells.Item(77, 8).Value = 54409.285
$ws.Cells.Item(77, 9).Value = 19999
$ws.Cells.Item(77, 10).Value = 80217
$ws.Cells.Item(77, 11).Value = 59997
$ws.Cells.Item(77, 12).Value = 240651
$ws.Cells.Item(77, 13).Value = -55005
$ws.Cells.Item(77, 14).Value = -250635

# Row 113
$ws.Cells.Item(113, 8).Value = 1668.909
$ws.Cells.Item(113, 9).Value = 1401.125
$ws.Cells.Item(113, 11).Value = 1401.125
$ws.Cells.Item(113, 13).Value = 768.875

$ws = $wb.Worksheets.Item("WVR")
# Row 126
$ws.Cells.Item(126, 8).Value = 2644.6667
$ws.Cells.Item(126, 9).Value = 3298.5
$ws.Cells.Item(126, 11).Value = 9895.5
$ws.Cells.Item(126, 13).Value = -7425.5

# Row 136
$ws.Cells.Item(136, 8).Value = 22001.5
$ws.Cells.Item(136, 9).Value = 15954
$ws.Cells.Item(136, 11).Value = 47862
$ws.Cells.Item(136, 13).Value = -45312
